# "Updated symbol list" run (Sun Dec 25 08:45:03 UTC 2022, GitHub Actions).
#
# The crypto-price sheet got re-scraped: a handful of Price (column D)
# values were refreshed, and the rows for the lowest-ranked coins shifted
# down by one slot (a new "One" entry was inserted at rank #8 / row 9,
# pushing WazirX -> MandalaExchangeToken -> ... -> CoinExToken each down
# one row, and similarly KickToken was inserted at rank #40 / row 41,
# pushing BKEXToken -> CEJI down one row each). Column D holds its prices
# as *text* (e.g. "245.18"), not numbers, so a plain Range.Value
# assignment of a numeric-looking string would silently be reinterpreted
# by Excel as a number (and would also tack on a "quote prefix" style the
# moment we force it back to text with a leading apostrophe). To keep
# these cells textual AND keep their original (unstyled) formatting, we
# assign with a leading "'" and then reset Style back to "Normal" right
# after, which drops the quote-prefix styling Excel otherwise applies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice($addr, $val) {
    # Force the numeric-looking string to stay plain text...
    $ws.Range($addr).Value = "'" + $val
    # ...without leaving behind the quote-prefix style Excel adds for that.
    $ws.Range($addr).Style = "Normal"
}

# Row 2 (BNB) - price refresh
Set-TextPrice "D2" "245.16"

# Row 3 (OKB) - price refresh
Set-TextPrice "D3" "23.05"

# Row 4 (HuobiToken) - price refresh
Set-TextPrice "D4" "5.408"

# Row 6 (GateToken) - price refresh
Set-TextPrice "D6" "3.391"

# Row 7 (MXToken) - price refresh
Set-TextPrice "D7" "0.8075"

# Row 8 (FTXToken) - price refresh
Set-TextPrice "D8" "0.9276"

# Row 9: new entry "One" inserted (was WazirX)
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextPrice "D9" "0.01123"
$ws.Range("E9").Value = "8OneONEBestin24h"

# Row 10: WazirX (was MandalaExchangeToken)
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextPrice "D10" "0.1422"
$ws.Range("E10").Value = "9WazirXWRX"

# Row 11: MandalaExchangeToken (was LiechtensteinCryptoassetsExchange)
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextPrice "D11" "0.07395"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

# Row 12: LiechtensteinCryptoassetsExchange (was BitrueCoin)
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextPrice "D12" "0.03378"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

# Row 13: BitrueCoin (was BitMartToken)
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextPrice "D13" "0.03038"
$ws.Range("E13").Value = "12BitrueCoinBTR"

# Row 14: BitMartToken (was MCDex)
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextPrice "D14" "0.09346"
$ws.Range("E14").Value = "13BitMartTokenBMX"

# Row 15: MCDex (was BitForexToken)
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextPrice "D15" "3.948"
$ws.Range("E15").Value = "14MCDexMCB"

# Row 16: BitForexToken (was CoinExToken)
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextPrice "D16" "0.001589"
$ws.Range("E16").Value = "15BitForexTokenBF"

# Row 17: CoinExToken (was One)
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextPrice "D17" "0.04805"
$ws.Range("E17").Value = "16CoinExTokenCET"

# Row 18 (TigerCash) - price refresh
Set-TextPrice "D18" "0.005339"

# Row 19 (HotbitToken) - price refresh
Set-TextPrice "D19" "0.004159"

# Row 20 (BitKan) - price refresh
Set-TextPrice "D20" "0.0009817"

# Row 22 (LEO) - price refresh
Set-TextPrice "D22" "3.659"

# Row 23 (KuCoinToken) - price refresh
Set-TextPrice "D23" "6.450"

# Row 24 (BTSEToken) - price refresh
Set-TextPrice "D24" "2.186"

# Row 40 (IDEX) - price refresh
Set-TextPrice "D40" "0.03950"

# Row 41: new entry "KickToken" inserted (was BKEXToken)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextPrice "D41" "0.006242"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42: BKEXToken (was CEJI)
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextPrice "D42" "0.1073"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43: CEJI (was KickToken)
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextPrice "D43" "0.002901"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44 (LocalTraders) - price refresh
Set-TextPrice "D44" "0.006939"

# Row 45 (CoinLion) - price refresh
Set-TextPrice "D45" "0.00005204"

# Row 48 (CoinbaseStockToken) - trailing "Bestin24h" badge removed from the volume label
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"

# Row 49 (BOLO) - price refresh
Set-TextPrice "D49" "0.002026"
